# Add 2022-Q4 data: new sheet right after "总计", plus a summary row in "总计".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new "2022-Q4" worksheet right after "总计" (first sheet).
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)
$q4Sheet = $wb.Worksheets.Add($null, $totalSheet)
$q4Sheet.Name = "2022-Q4"

# Header row (row 1), columns B:H - bold/centered/bordered like the other quarter sheets.
$q4Headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$q4Cols = @("B", "C", "D", "E", "F", "G", "H")
for ($i = 0; $i -lt $q4Headers.Length; $i++) {
    $q4Sheet.Range($q4Cols[$i] + "1").Value = $q4Headers[$i]
}

# Fund rows data: index, code, name, scale, position, ratio, marketValue, rank
$q4Data = @(
    @(0,  "002685", "中欧丰泓沪港深灵活配置混合A", "54.86", "93.29", "5.52", "3.0283", 4),
    @(1,  "002686", "中欧丰泓沪港深灵活配置混合C", "16.87", "93.29", "5.52", "0.9312", 4),
    @(2,  "013991", "中欧港股通精选一年持有混合A", "7.23",  "91.98", "5.84", "0.4222", 4),
    @(3,  "013992", "中欧港股通精选一年持有混合C", "4.91",  "91.98", "5.84", "0.2867", 4),
    @(4,  "671010", "西部利得策略优选混合A",       "1.88",  "92.90", "9.69", "0.1822", 3),
    @(5,  "011203", "永赢惠添益混合A",             "4.30",  "92.23", "3.69", "0.1587", 9),
    @(6,  "003413", "华泰柏瑞新经济沪港深混合",     "1.44",  "94.26", "7.80", "0.1123", 4),
    @(7,  "016570", "嘉实价值丰润混合A",           "2.11",  "63.89", "2.99", "0.0631", 10),
    @(8,  "011060", "西部利得策略优选混合C",       "0.48",  "92.90", "9.69", "0.0465", 3),
    @(9,  "011355", "华泰柏瑞港股通时代机遇混合A", "0.70",  "94.61", "6.52", "0.0456", 5),
    @(10, "011204", "永赢惠添益混合C",             "1.11",  "92.23", "3.69", "0.0410", 9),
    @(11, "460010", "华泰柏瑞亚洲领导企业混合（QDII）", "0.52", "97.17", "6.31", "0.0328", 2),
    @(12, "011356", "华泰柏瑞港股通时代机遇混合C", "0.39",  "94.61", "6.52", "0.0254", 5),
    @(13, "016298", "中欧丰泰港股通混合C",         "0.43",  "93.01", "5.70", "0.0245", 5),
    @(14, "016297", "中欧丰泰港股通混合A",         "0.27",  "93.01", "5.70", "0.0154", 5),
    @(15, "016571", "嘉实价值丰润混合C",           "0.20",  "63.89", "2.99", "0.0060", 10)
)

foreach ($row in $q4Data) {
    $r = [int]$row[0] + 2
    $q4Sheet.Range("A" + $r).Value = [int]$row[0]
    $q4Sheet.Range("B" + $r).Value = "'" + $row[1]
    $q4Sheet.Range("C" + $r).Value = $row[2]
    $q4Sheet.Range("D" + $r).Value = "'" + $row[3]
    $q4Sheet.Range("E" + $r).Value = "'" + $row[4]
    $q4Sheet.Range("F" + $r).Value = "'" + $row[5]
    $q4Sheet.Range("G" + $r).Value = "'" + $row[6]
    $q4Sheet.Range("H" + $r).Value = [int]$row[7]
}

# Match the look of the other quarter sheets: bold/centered/bordered header row
# and index column (copy formatting from the equivalent cells on "2022-Q3").
$q3Sheet = $wb.Worksheets.Item("2022-Q3")
$q3Sheet.Range("B1:H1").Copy()
$q4Sheet.Range("B1:H1").PasteSpecial(-4122)
$q3Sheet.Range("A2").Copy()
$q4Sheet.Range("A2:A17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2. Insert the new 2022-Q4 summary row into "总计" (row 2), shifting the
#    existing quarters down, and keep the running index (column A) sequential.
# ---------------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("B2:D2").ClearFormats()
$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("C2").Value = 16
$totalSheet.Range("D2").Value = 5.42

# Restore column-A's bold/centered/bordered look on the newly inserted row.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
